#
# githubProba.docx update
#
# The commit only makes one visible content change to the document body:
# the paragraph that reads
#     "proba TC 196 vjgjgjg ghfhfg nhg "
# gains extra trailing text and becomes
#     "proba TC 196 vjgjgjg ghfhfg nhg proabaaaaaaaaa "
#
# (Everything else in the raw OOXML diff -- the __RefHeading___Toc... bookmark
# names, the matching TOC hyperlink anchors, and the _Toc##############
# caption/heading bookmark ids -- are Word's own internal, auto-regenerated
# TOC/heading bookmark identifiers. They change any time Word recalculates
# the table of contents on save and carry no user-visible meaning; they are
# not something a document edit targets directly.)
#
$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "proba TC 196 vjgjgjg ghfhfg nhg ",  # FindText
    $true,                                # MatchCase
    $false,                               # MatchWholeWord
    $false,                               # MatchWildcards
    $false,                               # MatchSoundsLike
    $false,                               # MatchAllWordForms
    $true,                                # Forward
    1,                                    # Wrap (wdFindContinue)
    $false,                               # Format
    "proba TC 196 vjgjgjg ghfhfg nhg proabaaaaaaaaa ",  # ReplaceWith
    2                                     # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Target paragraph text was not found; document may differ from expected baseline."
}
